$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 84

# Copy formatting from the row above (row 83) into the new row 84,
# then set the values for the new match record.
$ws.Range("A83:V83").Copy() | Out-Null
$ws.Range("A84:V84").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 83
$ws.Cells.Item($row, 2).Value = "bosnia-and-herzegovina"
$ws.Cells.Item($row, 3).Value = "premijer-liga-bih"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45254.75
$ws.Cells.Item($row, 6).Value = "FK Sarajevo"
$ws.Cells.Item($row, 7).Value = 6
$ws.Cells.Item($row, 8).Value = "Tuzla City"
$ws.Cells.Item($row, 9).Value = 2
$ws.Cells.Item($row, 10).Value = 1.33
$ws.Cells.Item($row, 11).Value = "23/11/2023 07:12"
$ws.Cells.Item($row, 12).Value = 1.42
$ws.Cells.Item($row, 13).Value = "24/11/2023 17:56"
$ws.Cells.Item($row, 14).Value = 4.63
$ws.Cells.Item($row, 15).Value = "23/11/2023 07:12"
$ws.Cells.Item($row, 16).Value = 4.28
$ws.Cells.Item($row, 17).Value = "24/11/2023 17:57"
$ws.Cells.Item($row, 18).Value = 6.65
$ws.Cells.Item($row, 19).Value = "23/11/2023 07:12"
$ws.Cells.Item($row, 20).Value = 7.65
$ws.Cells.Item($row, 21).Value = "24/11/2023 17:57"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/fk-sarajevo-tuzla-city/CvIVXWxb/"
